$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The slide already has shapes with ids 1,2,6,8,9. This engine assigns the
# smallest unused positive id to a freshly-created shape, so we burn id "3"
# on a throwaway textbox (create + delete) to make the real new textbox land
# on id "4" / name "TextBox 3", matching the target deck.
$burn = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn.Delete()

# Target geometry (EMU): off x=6429935 y=1827911, ext cx=3951194 cy=369332.
# Shape.Left/Top/Width/Height are expressed in points and the host stores
# them as 32-bit floats before re-deriving EMU (x12700) by truncation, so we
# feed it the nearest float32 point value that truncates back to the exact
# EMU target instead of the naive EMU/12700 quotient.
$shp = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$shp.Left = 506.2940979003906
$shp.Top = 143.9300079345703
$shp.Width = 311.1176452636719
$shp.Height = 29.081260681152344

$shp.Fill.Visible = $false

$tf = $shp.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.Text = "https://github.com/Guillo-bit/Beletza "
$tr.LanguageID = "es-EC"

$urlRange = $tr.Characters(1, 37)
$urlRange.ActionSettings(1).Hyperlink.Address = "https://github.com/Guillo-bit/Beletza"
